# Update cryptos list (price + volume figures refreshed).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.284.32'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.846.73'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '241.27'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6735'
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07449'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2955'
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.95'
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07718'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.839.68'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.010'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6730'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '86.27'
$ws.Range('E15').Value = '  -1.46%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.153'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '29.278.21'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008340'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '228.95'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.54'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.211'
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '161.14'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.722'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E26').Value = '  -3.46%  '
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.512'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.071'
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05310'
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7598'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.877'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.139'
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.673'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.325.60'
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.726'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9184'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.988'
$ws.Range('E41').Value = '  +3.74%  '
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '103.38'
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.08158'
$ws.Range('E44').Value = '  +11.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.994.57'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5167'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.783'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '64.24'
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00000000121'
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.150'
$ws.Range('E50').Value = '  -3.79%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05957'
$ws.Range('E51').Value = '  +0.24%  '
